$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 values (mirrors the new <row r="6"> block added to sheet1.xml)
$ws.Range("A6").Value = 46007.44570207176
$ws.Range("A6").NumberFormat = $ws.Range("A5").NumberFormat

$ws.Range("B6").Value = "'"
$ws.Range("B6").Style = "Normal"

$ws.Range("C6").Value = "COMISARIA 9"

$ws.Range("D6").Value = "'54"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "REPARACIÓN"

$ws.Range("F6").Value = "BXCV"

$ws.Range("G6").Value = "TALLER POLICIAL"

$ws.Range("H6").Value = "'"
$ws.Range("H6").Style = "Normal"

$ws.Range("I6").Value = "INGRESADO"

$ws.Range("J6").Value = "'"
$ws.Range("J6").Style = "Normal"
